$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)

$notes = $s.NotesPage
$body = $notes.Shapes.Placeholders.Item(2)
$body.TextFrame.TextRange.Text = "=== Comparación estricta. Tipo y valor."
